$d = $word.ActiveDocument

# Paragraph 4: "Refer AssignmentQ1.html file" -> "Refer AssignmentQ1.html, AssignmentQ1.ts, AssignmentQ1.js file, "
$p = $d.Paragraphs.Item(4)
$p.Range.Find.Execute("Refer AssignmentQ1.html file", $true, $false, $false, $false, $false, $true, 1, $false, "Refer AssignmentQ1.html, AssignmentQ1.ts, AssignmentQ1.js file, ", 2)

# Paragraph 6: "Refer AssignmentQ2.html file" -> "Refer AssignmentQ2.html, AssignmentQ2.ts, AssignmentQ2.js file"
$p = $d.Paragraphs.Item(6)
$p.Range.Find.Execute(".html file", $true, $false, $false, $false, $false, $true, 1, $false, ".html, AssignmentQ2.ts, AssignmentQ2.js file", 2)

# Paragraph 8: "Refer AssignmentQ3.html file" -> "Refer AssignmentQ3.html, AssignmentQ3.ts, AssignmentQ3.js file"
$p = $d.Paragraphs.Item(8)
$p.Range.Find.Execute(".html file", $true, $false, $false, $false, $false, $true, 1, $false, ".html, AssignmentQ3.ts, AssignmentQ3.js file", 2)

# Paragraph 12: "Refer AssignmentQ4.html file" -> "Refer AssignmentQ4.html, AssignmentQ4.ts, AssignmentQ4.js file"
$p = $d.Paragraphs.Item(12)
$p.Range.Find.Execute(".html file", $true, $false, $false, $false, $false, $true, 1, $false, ".html, AssignmentQ4.ts, AssignmentQ4.js file", 2)

# Paragraph 17: "Refer AssignmentQ5.html file" -> "Refer AssignmentQ5.html, AssignmentQ5.ts, AssignmentQ5.js file"
$p = $d.Paragraphs.Item(17)
$p.Range.Find.Execute(".html file", $true, $false, $false, $false, $false, $true, 1, $false, ".html, AssignmentQ5.ts, AssignmentQ5.js file", 2)

# Paragraph 21: "Refer AssignmentQ6.html file" -> "Refer AssignmentQ6.html, AssignmentQ6.ts, AssignmentQ6.js file"
$p = $d.Paragraphs.Item(21)
$p.Range.Find.Execute(".html file", $true, $false, $false, $false, $false, $true, 1, $false, ".html, AssignmentQ6.ts, AssignmentQ6.js file", 2)

# Paragraph 26: "Refer AssignmentQ7.html file" -> "Refer AssignmentQ7.html, AssignmentQ7.ts, AssignmentQ7.js file"
$p = $d.Paragraphs.Item(26)
$p.Range.Find.Execute("Refer AssignmentQ7.html file", $true, $false, $false, $false, $false, $true, 1, $false, "Refer AssignmentQ7.html, AssignmentQ7.ts, AssignmentQ7.js file", 2)

# Paragraph 29: "Refer AssignmentQ8.html file" -> "Refer AssignmentQ8.html, AssignmentQ8.ts, AssignmentQ8.js file"
$p = $d.Paragraphs.Item(29)
$p.Range.Find.Execute(".html file", $true, $false, $false, $false, $false, $true, 1, $false, ".html, AssignmentQ8.ts, AssignmentQ8.js file", 2)

Write-Host "Done"
